# HC Services.xlsx - "Add files via upload" edit
#
# Applies the data / formatting changes described by the target diff:
#   - Main!D6 gets a new value (254.03)
#   - Main!A14, A16, A18, A20 lose their "x" marker; A22, A23, A24, A25, A26,
#     A27, A30 gain one
#   - Main column widths B/C nudged slightly
#   - Main's frozen-pane selection moves from J7 to H16
#   - CDMO!B2/C2 get header labels ("Name"/"Ticker"); A5/A6 lose their "x"
#     marker
#   - CDMO gains column widths (A=5, B=17) and is zoomed to 130%
#   - Resources / Private / CDMO sheets lose their stray saved cell
#     selections (reset to A1)
#   - Main stays the active / selected tab throughout

$wb = $excel.ActiveWorkbook

$wsMain = $wb.Worksheets.Item("Main")
$wsResources = $wb.Worksheets.Item("Resources")
$wsPrivate = $wb.Worksheets.Item("Private")
$wsCdmo = $wb.Worksheets.Item("CDMO")

# ---------------------------------------------------------------------------
# Main: new data point in row 6
# ---------------------------------------------------------------------------
$wsMain.Range("D6").Value = 254.03

# ---------------------------------------------------------------------------
# Main: the "x" flag column (A) moves off rows 14/16/18/20 and onto
# 22/23/24/25/26/27/30
# ---------------------------------------------------------------------------
$wsMain.Range("A14").ClearContents()
$wsMain.Range("A16").ClearContents()
$wsMain.Range("A18").ClearContents()
$wsMain.Range("A20").ClearContents()

$wsMain.Range("A22").Value = "x"
$wsMain.Range("A23").Value = "x"
$wsMain.Range("A24").Value = "x"
$wsMain.Range("A25").Value = "x"
$wsMain.Range("A26").Value = "x"
$wsMain.Range("A27").Value = "x"
$wsMain.Range("A30").Value = "x"

# ---------------------------------------------------------------------------
# Main: small column-width nudges (B and C)
# ---------------------------------------------------------------------------
$wsMain.Columns.Item(2).ColumnWidth = 18
$wsMain.Columns.Item(3).ColumnWidth = 9.666666666666666

# ---------------------------------------------------------------------------
# CDMO: add header labels in row 2, drop the "x" flags on rows 5/6
# ---------------------------------------------------------------------------
$wsCdmo.Range("B2").Value = "Name"
$wsCdmo.Range("C2").Value = "Ticker"
$wsCdmo.Range("A5").ClearContents()
$wsCdmo.Range("A6").ClearContents()

# ---------------------------------------------------------------------------
# CDMO: new column widths + 130% zoom, then reset its lingering selection
# ---------------------------------------------------------------------------
$wsCdmo.Columns.Item(1).ColumnWidth = 4.166666666666667
$wsCdmo.Columns.Item(2).ColumnWidth = 16.166666666666668

$wsCdmo.Activate()
$excel.ActiveWindow.Zoom = 130
$wsCdmo.Range("A1").Select()

# ---------------------------------------------------------------------------
# Resources / Private: clear the stray saved selections
# ---------------------------------------------------------------------------
$wsResources.Activate()
$wsResources.Range("A1").Select()

$wsPrivate.Activate()
$wsPrivate.Range("A1").Select()

# ---------------------------------------------------------------------------
# Leave Main as the active sheet/tab, with the new frozen-pane selection
# ---------------------------------------------------------------------------
$wsMain.Activate()
$wsMain.Range("H16").Select()
